$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 652.25
$ws.Range("I41").Value = 280.4
$ws.Range("J41").Value = 917.8570999999999
$ws.Range("K41").Value = 280.4
$ws.Range("L41").Value = 917.8570999999999
$ws.Range("M41").Value = 159.6
$ws.Range("N41").Value = -1797.8571

$ws.Range("H64").Value = 4677.5713
$ws.Range("I64").Value = 2148
$ws.Range("J64").Value = 6574.75
$ws.Range("K64").Value = 2148
$ws.Range("L64").Value = 6574.75
$ws.Range("M64").Value = -1900
$ws.Range("N64").Value = -7070.75

$ws.Range("H67").Value = 4677.5713
$ws.Range("I67").Value = 2148
$ws.Range("J67").Value = 6574.75
$ws.Range("K67").Value = 2148
$ws.Range("L67").Value = 6574.75
$ws.Range("M67").Value = -1290
$ws.Range("N67").Value = -8290.75

$ws.Range("H106").Value = 3726.2727
$ws.Range("J106").Value = 2833
$ws.Range("L106").Value = 2833
$ws.Range("N106").Value = -4095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1765579.5
$ws.Range("I61").Value = 4262.5
$ws.Range("K61").Value = 4262.5
$ws.Range("M61").Value = -4050.5

$ws.Range("H74").Value = 24212.137
$ws.Range("I74").Value = 1846.3334
$ws.Range("J74").Value = 51051.1
$ws.Range("K74").Value = 1846.3334
$ws.Range("L74").Value = 51051.1
$ws.Range("M74").Value = -972.3334
$ws.Range("N74").Value = -52799.1

$ws.Range("H77").Value = 24212.137
$ws.Range("I77").Value = 1846.3334
$ws.Range("J77").Value = 51051.1
$ws.Range("K77").Value = 9231.666999999999
$ws.Range("L77").Value = 255255.5
$ws.Range("M77").Value = -4863.666999999999
$ws.Range("N77").Value = -263991.5

$ws.Range("H97").Value = 1033.421
$ws.Range("I97").Value = 821.5833
$ws.Range("K97").Value = 821.5833
$ws.Range("M97").Value = -325.5833

$ws.Range("H122").Value = 3210.7693
$ws.Range("I122").Value = 1586.2941
$ws.Range("J122").Value = 6279.222
$ws.Range("K122").Value = 4758.8823
$ws.Range("L122").Value = 18837.666
$ws.Range("M122").Value = -2308.8823
$ws.Range("N122").Value = -23737.666

$ws.Range("H136").Value = 1765579.5
$ws.Range("I136").Value = 4262.5
$ws.Range("K136").Value = 12787.5
$ws.Range("M136").Value = -10237.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1471.0667
$ws.Range("I105").Value = 1347.25
$ws.Range("K105").Value = 1347.25
$ws.Range("M105").Value = 399.75

$ws.Range("H107").Value = 1421.56
$ws.Range("I107").Value = 1335.4375
$ws.Range("J107").Value = 1574.6666
$ws.Range("K107").Value = 1335.4375
$ws.Range("L107").Value = 1574.6666
$ws.Range("M107").Value = 584.5625
$ws.Range("N107").Value = -5414.6666

$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 33209152
$ws.Range("I132").Value = 2007.238
$ws.Range("K132").Value = 6021.714
$ws.Range("M132").Value = -3491.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7144.909
$ws.Range("I70").Value = 4726.7144
$ws.Range("K70").Value = 4726.7144
$ws.Range("M70").Value = -4456.7144

$ws.Range("H73").Value = 7144.909
$ws.Range("I73").Value = 4726.7144
$ws.Range("K73").Value = 4726.7144
$ws.Range("M73").Value = -3790.7144

$ws.Range("H97").Value = 1875.375
$ws.Range("I97").Value = 1769.4375
$ws.Range("J97").Value = 2087.25
$ws.Range("K97").Value = 1769.4375
$ws.Range("L97").Value = 2087.25
$ws.Range("M97").Value = -1273.4375
$ws.Range("N97").Value = -3079.25

$ws.Range("H126").Value = 7673.1904
$ws.Range("I126").Value = 10333.857
$ws.Range("J126").Value = 6342.857
$ws.Range("K126").Value = 31001.571
$ws.Range("L126").Value = 19028.571
$ws.Range("M126").Value = -28531.571
$ws.Range("N126").Value = -23968.571

$ws.Range("H132").Value = 693279.8
$ws.Range("I132").Value = 6271.0625
$ws.Range("J132").Value = 1914628.8
$ws.Range("K132").Value = 18813.1875
$ws.Range("L132").Value = 5743886.4
$ws.Range("M132").Value = -16283.1875
$ws.Range("N132").Value = -5748946.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5141.4136
$ws.Range("I7").Value = 4876.476
$ws.Range("K7").Value = 4876.476
$ws.Range("M7").Value = -4764.476

$ws.Range("H16").Value = 2835.75
$ws.Range("J16").Value = 4084
$ws.Range("L16").Value = 4084
$ws.Range("N16").Value = -4424

$ws.Range("H22").Value = 5349.636
$ws.Range("I22").Value = 6506.2666
$ws.Range("K22").Value = 6506.2666
$ws.Range("M22").Value = -6211.2666

$ws.Range("H27").Value = 5349.636
$ws.Range("I27").Value = 6506.2666
$ws.Range("K27").Value = 6506.2666
$ws.Range("M27").Value = -6399.2666

$ws.Range("H40").Value = 3752.25
$ws.Range("I40").Value = 3752.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3752.25
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3616.25
$ws.Range("N40").ClearContents()

$ws.Range("H53").Value = 15000
$ws.Range("I53").Value = 12000
$ws.Range("K53").Value = 12000
$ws.Range("M53").Value = -11482

$ws.Range("H64").Value = 300000
$ws.Range("J64").Value = 300000
$ws.Range("L64").Value = 300000
$ws.Range("N64").Value = -300450

$ws.Range("H67").Value = 300000
$ws.Range("J67").Value = 300000
$ws.Range("L67").Value = 300000
$ws.Range("N67").Value = -301560

$ws.Range("H100").Value = 2406.6191
$ws.Range("J100").Value = 2653.4285
$ws.Range("L100").Value = 2653.4285
$ws.Range("N100").Value = -3735.4285

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H122").Value = 4844.943
$ws.Range("I122").Value = 3480.5293
$ws.Range("J122").Value = 6133.5557
$ws.Range("K122").Value = 10441.5879
$ws.Range("L122").Value = 18400.6671
$ws.Range("M122").Value = -7991.5879
$ws.Range("N122").Value = -23300.6671

$ws.Range("H126").Value = 5141.4136
$ws.Range("I126").Value = 4876.476
$ws.Range("K126").Value = 14629.428
$ws.Range("M126").Value = -12159.428

$ws.Range("H132").Value = 2589524.8
$ws.Range("J132").Value = 9976051
$ws.Range("L132").Value = 29928153
$ws.Range("N132").Value = -29933213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16258.571
$ws.Range("J41").Value = 15574
$ws.Range("L41").Value = 15574
$ws.Range("N41").Value = -16354

$ws.Range("H62").Value = 26999.875
$ws.Range("J62").Value = 26250
$ws.Range("L62").Value = 26250
$ws.Range("N62").Value = -27498

$ws.Range("H65").Value = 26999.875
$ws.Range("J65").Value = 26250
$ws.Range("L65").Value = 131250
$ws.Range("N65").Value = -137490

$ws.Range("H132").Value = 611881.75
$ws.Range("I132").Value = 2800.4
$ws.Range("K132").Value = 8401.200000000001
$ws.Range("M132").Value = -5871.200000000001

$ws.Range("H136").Value = 876876.2
$ws.Range("I136").Value = 2923.1667
$ws.Range("K136").Value = 8769.500100000001
$ws.Range("M136").Value = -6219.500100000001
